# EscolherConfiguracaoOtima.xlsx
#
# Commit message: "Alterei alt para opt e para break"
# (Changed the UML "alt" fragment label to an "opt"-style exception label,
#  and turned the old "Regressa a 1" return-flow note into an empty cell,
#  i.e. a "break" with no further text.)
#
# Concretely, in the use-case description table:
#   - B20 ("Alternativo 2 [Não confirma compra] (passo 11)")
#         becomes "Exceção 1 [Não confirma compra] (passo 11)"
#   - C20 ("11.1 Não confirma compra") stays the same text
#   - D21 ("Regressa a 1") is cleared out entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the alternative-fragment heading in B20.
$ws.Range("B20").Value = "Exceção 1 [Não confirma compra] (passo 11)"

# Remove the old "Regressa a 1" note, leaving the cell blank (formatting kept).
$ws.Range("D21").ClearContents()

# Reflect the author's final selection position in the sheet view.
$ws.Range("B26").Select()
